$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194, shifting existing rows 194-249 down to 195-250.
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new record.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R match the record that used to be
# (and still is, one row below) at row 195, i.e. the same market/category metadata.
$ws.Cells.Item(194, 1).Value = 5
$ws.Cells.Item(194, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(194, 3).Value = "Maule"
$ws.Cells.Item(194, 4).Value = 44841
$ws.Cells.Item(194, 5).Value = 7
$ws.Cells.Item(194, 6).Value = 100112017
$ws.Cells.Item(194, 7).Value = "Apio"
$ws.Cells.Item(194, 8).Value = "Americana (o)"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 700
$ws.Cells.Item(194, 11).Value = 9000
$ws.Cells.Item(194, 12).Value = 9000
$ws.Cells.Item(194, 13).Value = 9000
$ws.Cells.Item(194, 14).Value = "$/docena de matas"
$ws.Cells.Item(194, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(194, 16).Value = 1500
$ws.Cells.Item(194, 17).Value = 6
$ws.Cells.Item(194, 18).Value = "Hortaliza"

# Make sure the date cell keeps the existing date-style formatting (style index 2
# used throughout column D) now that a fresh cell was inserted.
$ws.Cells.Item(194, 4).NumberFormat = $ws.Cells.Item(195, 4).NumberFormat
